$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "59.699.41"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.651.28"
$ws.Range("E3").Value = "  +1.71%  "
Set-TextValue "D5" "537.51"
$ws.Range("E5").Value = "  -1.33%  "
Set-TextValue "D6" "146.47"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +1.27%  "
Set-TextValue "D9" "6.83"
$ws.Range("E9").Value = "  +5.65%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "3.123.91"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "59.612.55"
$ws.Range("E14").Value = "  +0.62%  "
Set-TextValue "D15" "21.44"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").Value = "2.678.13"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("E17").Value = "  +1.09%  "
Set-TextValue "D18" "4.48"
$ws.Range("E18").Value = "  +2.82%  "
Set-TextValue "D19" "340.10"
$ws.Range("E19").Value = "  -0.87%  "
Set-TextValue "D20" "10.36"
$ws.Range("E20").Value = "  +2.27%  "
Set-TextValue "D21" "6.22"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("E22").Value = "  +0.02%  "
Set-TextValue "D23" "66.66"
$ws.Range("E23").Value = "  -1.30%  "
Set-TextValue "D24" "0.418"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("E25").Value = "  -0.47%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  -0.09%  "
Set-TextValue "D27" "7.31"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  +1.76%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -2.94%  "
Set-TextValue "D31" "5.88"
$ws.Range("E31").Value = "  +1.16%  "
Set-TextValue "D32" "18.90"
$ws.Range("E32").Value = "  +0.72%  "
Set-TextValue "D33" "150.92"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  +2.80%  "
Set-TextValue "D36" "0.842"
$ws.Range("E36").Value = "  +3.56%  "
Set-TextValue "D37" "0.842"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -0.90%  "
Set-TextValue "D39" "3.62"
$ws.Range("E39").Value = "  +1.82%  "
Set-TextValue "D40" "286.95"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("E41").Value = "  -0.09%  "
Set-TextValue "D42" "0.608"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Value = "1.967.21"
$ws.Range("E48").Value = "  +1.16%  "
Set-TextValue "D49" "4.57"
$ws.Range("E49").Value = "  +1.39%  "
Set-TextValue "D50" "18.43"
$ws.Range("E50").Value = "  +0.37%  "
Set-TextValue "D51" "112.24"
$ws.Range("E51").Value = "  +1.26%  "
